$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures for the existing last quarter (01-01-2021), row 75
$ws.Range("I75").Value = 1640
$ws.Range("J75").Value = 162
$ws.Range("L75").Value = 9589

# Append the new quarter (01-04-2021) as row 76.
# A76 must hold the literal text "01-04-2021" (not an auto-converted date).
# Entering it as a formula that evaluates to the text, then pasting the
# result back as a value, keeps it as plain text without leaving any new
# number-format/style behind (unlike setting NumberFormat="@" directly).
$ws.Range("A76").Formula = '="01-04-2021"'
$ws.Range("A76").Copy()
$ws.Range("A76").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B76").Value = 2030
$ws.Range("C76").Value = 46
$ws.Range("D76").Value = 1983
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 26038
$ws.Range("G76").Value = 600
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 1529
$ws.Range("J76").Value = 137
$ws.Range("K76").Value = 330
$ws.Range("L76").Value = 9064
$ws.Range("M76").Value = 13207
$ws.Range("N76").Value = 1172
